# Adds a new inventory-shortage row ("NEVILOB 2.5 MG 14 TAB.") right above the
# existing "ROYAL REGIME ..." row (new physical row 21), pushing every row
# below it down by one, renumbering the running "#" column, and refreshing
# the grand-total cell accordingly. This mirrors a freshly regenerated
# "DaySale" export that now contains one additional line item.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a blank row at 21 - shifts old rows 21-33 down to 22-34,
#    carrying their values/styles/merges with them.
# ---------------------------------------------------------------------
$ws.Rows("21:21").Insert()

# The Insert above does not create merged ranges for the brand-new row,
# so recreate the usual A:B / C:G / H:K / L:M / N:O pattern used by every
# other data row.
$ws.Range("A21:B21").Merge()
$ws.Range("C21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("N21:O21").Merge()

# Copy the formatting (number formats, borders, fill, font, alignment) from
# the row directly beneath - which still carries the original data-row
# style - onto the new row, without touching its (still empty) values.
$ws.Range("A22:Q22").Copy()
$ws.Range("A21:Q21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Populate the new row with the NEVILOB item.
#    Columns A (plain number) and C/H/N/Q (already text-formatted cells)
#    can be written to directly. Columns L and P carry a *numeric* display
#    format even though their content is plain text ("1", "23.0000"), so a
#    direct .Value assignment would get auto-coerced to a real number.
#    Stage the text in a scratch cell (forced to Text via the leading
#    apostrophe) and paste-special just the *value* across, which keeps
#    the destination cell's existing number format/style untouched.
# ---------------------------------------------------------------------
$ws.Range("A21").Value = 15
$ws.Range("C21").Value = "NEVILOB 2.5 MG 14 TAB."
$ws.Range("H21").Value = "1:0"
$ws.Range("N21").Value = "46.00"
$ws.Range("Q21").Value = "0:1"

$scratch = $ws.Range("ZZ1:ZZ2")
$scratch.Value = "'1"
$ws.Range("ZZ1").Value = "'1"
$ws.Range("ZZ2").Value = "'23.0000"

$ws.Range("ZZ1").Copy()
$ws.Range("L21").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("ZZ2").Copy()
$ws.Range("P21").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0
$ws.Range("ZZ1:ZZ2").ClearContents()

# ---------------------------------------------------------------------
# 3) Renumber the running "#" column for every row that shifted down
#    (old #15 .. #25 become #16 .. #26).
# ---------------------------------------------------------------------
for ($r = 22; $r -le 32; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 1).Value + 1
}

# ---------------------------------------------------------------------
# 4) Refresh the grand total (was 1011.47, now +23.00 for the new row).
# ---------------------------------------------------------------------
$ws.Range("P33").Value = 1034.47
